# Generate Report for Handback
# Update the handback timestamps for the "0fea17c1-a2ec-43a1-9a95-dbaf8bd228db.md"
# file's row (row 2) across the Overview, zh-cn and de-de sheets, reflecting a
# fresh handback xliff generation.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-06 00:54:31"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-06 00:54:26"
$wsZhCn.Range("K2").Value = "2016-09-06 00:54:44"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-06 00:54:31"
$wsDeDe.Range("K2").Value = "2016-09-06 00:54:52"
